$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3855
$ws.Range("E2").Value = 404
$ws.Range("F2").Value = 404
$ws.Range("G2").Value = 462
$ws.Range("H2").Value = 354
$ws.Range("I2").Value = 334
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 4374
$ws.Range("L2").Value = 2494
$ws.Range("M2").Value = 1880
$ws.Range("N2").Value = 1803
$ws.Range("O2").Value = 77
$ws.Range("P2").Value = 58
$ws.Range("Q2").Value = 1037
$ws.Range("R2").Value = -395
$ws.Range("S2").Value = -317
$ws.Range("T2").Value = 103
$ws.Range("U2").Value = 933
$ws.Range("V2").Value = 193
$ws.Range("W2").Value = 10.49
$ws.Range("X2").Value = 9.18
$ws.Range("Y2").Value = 18.25
$ws.Range("Z2").Value = 8.699999999999999
$ws.Range("AA2").Value = 132.66
$ws.Range("AB2").Value = 3628.43
$ws.Range("AC2").Value = 2879
$ws.Range("AD2").Value = 27.06
$ws.Range("AE2").Value = 16369
$ws.Range("AF2").Value = 4.76
$ws.Range("AG2").Value = 1300
$ws.Range("AH2").Value = 1.67
$ws.Range("AI2").Value = 43.3
$ws.Range("AJ2").Value = 11616185

# Row 3
$ws.Range("D3").Value = 4594
$ws.Range("E3").Value = 448
$ws.Range("F3").Value = 448
$ws.Range("G3").Value = 469
$ws.Range("H3").Value = 342
$ws.Range("I3").Value = 317
$ws.Range("J3").Value = 24
$ws.Range("K3").Value = 5044
$ws.Range("L3").Value = 2871
$ws.Range("M3").Value = 2173
$ws.Range("N3").Value = 1969
$ws.Range("O3").Value = 204
$ws.Range("P3").Value = 58
$ws.Range("Q3").Value = 350
$ws.Range("R3").Value = -9
$ws.Range("S3").Value = -24
$ws.Range("T3").Value = 164
$ws.Range("U3").Value = 186
$ws.Range("V3").Value = 305
$ws.Range("W3").Value = 9.74
$ws.Range("X3").Value = 7.44
$ws.Range("Y3").Value = 16.83
$ws.Range("Z3").Value = 7.26
$ws.Range("AA3").Value = 132.11
$ws.Range("AB3").Value = 3889.78
$ws.Range("AC3").Value = 2732
$ws.Range("AD3").Value = 42.09
$ws.Range("AE3").Value = 17837
$ws.Range("AF3").Value = 6.45
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 1.3
$ws.Range("AI3").Value = 52.17
$ws.Range("AJ3").Value = 11616185

# Row 4
$ws.Range("D4").Value = 5955
$ws.Range("E4").Value = 209
$ws.Range("F4").Value = 209
$ws.Range("G4").Value = 239
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 80
$ws.Range("J4").Value = -5
$ws.Range("K4").Value = 5202
$ws.Range("L4").Value = 3082
$ws.Range("M4").Value = 2121
$ws.Range("N4").Value = 1895
$ws.Range("O4").Value = 226
$ws.Range("P4").Value = 58
$ws.Range("Q4").Value = 275
$ws.Range("R4").Value = 67
$ws.Range("S4").Value = -255
$ws.Range("T4").Value = 480
$ws.Range("U4").Value = -205
$ws.Range("V4").Value = 202
$ws.Range("W4").Value = 3.52
$ws.Range("X4").Value = 1.26
$ws.Range("Y4").Value = 4.16
$ws.Range("Z4").Value = 1.47
$ws.Range("AA4").Value = 145.31
$ws.Range("AB4").Value = 3760.97
$ws.Range("AC4").Value = 692
$ws.Range("AD4").Value = 95.55
$ws.Range("AE4").Value = 17122
$ws.Range("AF4").Value = 3.86
$ws.Range("AG4").Value = 1500
$ws.Range("AH4").Value = 2.27
$ws.Range("AI4").Value = 206.58
$ws.Range("AJ4").Value = 11616185

# Row 5
$ws.Range("D5").Value = 8043
$ws.Range("E5").Value = 411
$ws.Range("F5").Value = 411
$ws.Range("G5").Value = 383
$ws.Range("H5").Value = 129
$ws.Range("I5").Value = 132
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 6282
$ws.Range("L5").Value = 3824
$ws.Range("M5").Value = 2457
$ws.Range("N5").Value = 2089
$ws.Range("O5").Value = 367
$ws.Range("P5").Value = 58
$ws.Range("Q5").Value = 662
$ws.Range("R5").Value = -218
$ws.Range("S5").Value = 329
$ws.Range("T5").Value = 214
$ws.Range("U5").Value = 448
$ws.Range("V5").Value = 255
$ws.Range("W5").Value = 5.11
$ws.Range("X5").Value = 1.61
$ws.Range("Y5").Value = 6.62
$ws.Range("Z5").Value = 2.25
$ws.Range("AA5").Value = 155.63
$ws.Range("AB5").Value = 4156.99
$ws.Range("AC5").Value = 1135
$ws.Range("AD5").Value = 90.75
$ws.Range("AE5").Value = 18876
$ws.Range("AF5").Value = 5.46
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 1.46
$ws.Range("AI5").Value = 125.91
$ws.Range("AJ5").Value = 11616185

# Row 6
$ws.Range("D6").Value = 8283
$ws.Range("E6").Value = 249
$ws.Range("F6").Value = 249
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 106
$ws.Range("I6").Value = 87
$ws.Range("K6").Value = 7161
$ws.Range("L6").Value = 4762
$ws.Range("M6").Value = 2399
$ws.Range("N6").Value = 1999
$ws.Range("P6").Value = 58
$ws.Range("Q6").Value = 202
$ws.Range("R6").Value = -1310
$ws.Range("S6").Value = 369
$ws.Range("T6").Value = 129
$ws.Range("U6").Value = 73
$ws.Range("V6").Value = 1224
$ws.Range("W6").Value = 3
$ws.Range("X6").Value = 1.28
$ws.Range("Y6").Value = 4.28
$ws.Range("Z6").Value = 1.57
$ws.Range("AA6").Value = 198.53
$ws.Range("AB6").Value = 3952.43
$ws.Range("AC6").Value = 753
$ws.Range("AD6").Value = 91.34
$ws.Range("AE6").Value = 18060
$ws.Range("AF6").Value = 3.81
$ws.Range("AG6").Value = 1300
$ws.Range("AH6").Value = 1.89
$ws.Range("AI6").Value = 164.43
$ws.Range("AJ6").Value = 11616185

# Row 7
$ws.Range("D7").Value = 7750
$ws.Range("E7").Value = 113
$ws.Range("G7").Value = -41
$ws.Range("H7").Value = -66
$ws.Range("I7").Value = -34
$ws.Range("K7").Value = 8417
$ws.Range("L7").Value = 6191
$ws.Range("M7").Value = 2227
$ws.Range("N7").Value = 1905
$ws.Range("P7").Value = 59
$ws.Range("Q7").Value = -87
$ws.Range("R7").Value = -431
$ws.Range("S7").Value = -173
$ws.Range("T7").Value = 693
$ws.Range("U7").Value = 67
$ws.Range("W7").Value = 1.46
$ws.Range("X7").Value = -0.85
$ws.Range("Y7").Value = -1.73
$ws.Range("Z7").Value = -0.85
$ws.Range("AA7").Value = 278.04
$ws.Range("AC7").Value = -290
$ws.Range("AD7").Value = -154.92
$ws.Range("AE7").Value = 17215
$ws.Range("AF7").Value = 2.61
$ws.Range("AG7").Value = 1144
$ws.Range("AH7").Value = 2.55
$ws.Range("AI7").Value = -394.87

# Row 8
$ws.Range("D8").Value = 8130
$ws.Range("E8").Value = 309
$ws.Range("G8").Value = 273
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 180
$ws.Range("K8").Value = 8937
$ws.Range("L8").Value = 6303
$ws.Range("M8").Value = 2634
$ws.Range("N8").Value = 2306
$ws.Range("P8").Value = 62
$ws.Range("Q8").Value = 702
$ws.Range("R8").Value = -376
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = 271
$ws.Range("U8").Value = 200
$ws.Range("W8").Value = 3.8
$ws.Range("X8").Value = 2.46
$ws.Range("Y8").Value = 8.57
$ws.Range("Z8").Value = 2.31
$ws.Range("AA8").Value = 239.32
$ws.Range("AC8").Value = 1552
$ws.Range("AD8").Value = 28.92
$ws.Range("AE8").Value = 20832
$ws.Range("AF8").Value = 2.16
$ws.Range("AG8").Value = 1200
$ws.Range("AH8").Value = 2.67
$ws.Range("AI8").Value = 77.3

# Row 9
$ws.Range("D9").Value = 8617
$ws.Range("E9").Value = 450
$ws.Range("G9").Value = 418
$ws.Range("H9").Value = 307
$ws.Range("I9").Value = 274
$ws.Range("K9").Value = 9021
$ws.Range("L9").Value = 6167
$ws.Range("M9").Value = 2853
$ws.Range("N9").Value = 2455
$ws.Range("P9").Value = 62
$ws.Range("Q9").Value = 718
$ws.Range("R9").Value = -343
$ws.Range("S9").Value = -227
$ws.Range("T9").Value = 253
$ws.Range("U9").Value = 210
$ws.Range("W9").Value = 5.23
$ws.Range("X9").Value = 3.57
$ws.Range("Y9").Value = 11.53
$ws.Range("Z9").Value = 3.42
$ws.Range("AA9").Value = 216.14
$ws.Range("AC9").Value = 2362
$ws.Range("AD9").Value = 19.01
$ws.Range("AE9").Value = 22187
$ws.Range("AF9").Value = 2.02
$ws.Range("AG9").Value = 1275
$ws.Range("AH9").Value = 2.84
$ws.Range("AI9").Value = 53.98

Write-Output "Done applying 244 cell updates"
